$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 244
$ws1.Range("F6").Value = 1167
$ws1.Range("F7").Value = 955
$ws1.Range("F10").Value = 81
$ws1.Range("F15").Value = 1384
$ws1.Range("F17").Value = 1286
$ws1.Range("F19").Value = 329
$ws1.Range("F20").Value = 1570
$ws1.Range("F26").Value = 1082
$ws1.Range("F27").Value = 378
$ws1.Range("F28").Value = 3349
$ws1.Range("F30").Value = 551

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 4

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 244
$ws4.Range("F10").Value = 1167
$ws4.Range("F11").Value = 955
$ws4.Range("F20").Value = 81
$ws4.Range("F27").Value = 1384
$ws4.Range("F29").Value = 1286
$ws4.Range("F31").Value = 329
$ws4.Range("F32").Value = 1570
$ws4.Range("F40").Value = 1082
$ws4.Range("F41").Value = 378
$ws4.Range("F42").Value = 3349
$ws4.Range("F44").Value = 551
$ws4.Range("F46").Value = 4
